# InvestmentCalc.xlsx update: revised investment assumptions
# (lower initial investment, higher revenue, lower costs, etc.)
# and the recalculated cash-flow figures that follow from them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Initial Investment ---
$ws.Range("B2").Value = -1500000

# --- Row 3: Depreciation (Year 1-10) ---
$ws.Range("C3:L3").Value = 45000

# --- Row 4: Incoming Payments (Year 1-10) ---
$ws.Range("C4:L4").Value = 700000

# --- Row 5: Outgoing Payments (Year 0 only changes) ---
$ws.Range("B5").Value = -70000

# --- Row 6: Residual (Year 10) ---
$ws.Range("L6").Value = 140000

# --- Row 7: restricted Equity ---
$ws.Range("B7").Value = -300000
$ws.Range("L7").Value = 300000

# --- Row 8: Yearly Net ---
$ws.Range("B8").Value = -1870000
$ws.Range("C8:K8").Value = 605000
$ws.Range("L8").Value = 1045000

# --- Row 9: Present Value ---
$ws.Range("B9").Value = -1870000
$ws.Range("C9").Value = 558118.0811808117
$ws.Range("D9").Value = 514869.0785800847
$ws.Range("E9").Value = 474971.4747048751
$ws.Range("F9").Value = 438165.5670709179
$ws.Range("G9").Value = 404211.7777406992
$ws.Range("H9").Value = 372889.0938567335
$ws.Range("I9").Value = 343993.6290191268
$ws.Range("J9").Value = 317337.296143106
$ws.Range("K9").Value = 292746.5831578468
$ws.Range("L9").Value = 466469.7316336035

# --- Row 10: Accumulated Present Value ---
$ws.Range("B10").Value = -1870000
$ws.Range("C10").Value = -1311881.918819188
$ws.Range("D10").Value = -797012.8402391034
$ws.Range("E10").Value = -322041.3655342283
$ws.Range("F10").Value = 116124.2015366896
$ws.Range("G10").Value = 520335.9792773888
$ws.Range("H10").Value = 893225.0731341223
$ws.Range("I10").Value = 1237218.702153249
$ws.Range("J10").Value = 1554555.998296355
$ws.Range("K10").Value = 1847302.581454202
$ws.Range("L10").Value = 2313772.313087805

# --- Row 11: Net Present Value ---
$ws.Range("B11").Value = 2313772.313087805

# --- Formatting: the Accumulated Present Value turns positive from year
#     4 onward (F10:L10), and the final NPV (B11) is now positive too, so
#     those cells switch from the "negative" red fill to the "positive"
#     green fill used elsewhere on the sheet. Copy the already-correct
#     formatting (number format + fill) from sibling cells that use each
#     target style, instead of hand-building a new style.
$ws.Range("C9").Copy()
$ws.Range("F10:L10").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("C8").Copy()
$ws.Range("B11").PasteSpecial(-4122)       # xlPasteFormats

$excel.CutCopyMode = 0
